$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$html = @"
<html lang="en">
 <head>
  <meta charset="utf-8"/>
  <meta content="width=device-width, initial-scale=1.0" name="viewport"/>
  <meta content="ie=edge" http-equiv="X-UA-Compatible"/>
  <title>
   Netflix
  </title>
  <link href="style.css" rel="stylesheet"/>
  <link href="mediaquery.css" rel="stylesheet"/>
  <link href="https://maxst.icons8.com/vue-static/landings/line-awesome/line-awesome/1.3.0/css/line-awesome.min.css" rel="stylesheet"/>
  <script crossorigin="anonymous" src="https://kit.fontawesome.com/bc3a1796c2.js">
  </script>
  <link href="https://image.flaticon.com/icons/png/512/870/870910.ico" rel="shortcut icon">
   <link href="https://cdnjs.cloudflare.com/ajax/libs/normalize/8.0.1/normalize.css" rel="stylesheet">
   </link>
  </link>
 </head>
 <body>
  <div class="navbar">
   <li class="logo">
    <img src="https://www.edigitalagency.com.au/wp-content/uploads/Netflix-logo-red-black-png.png"/>
   </li>
   <li class="buttons">
    Sign In
   </li>
  </div>
  <div class="main">
   <div class="area">
    <h1>
    </h1>
    <h3>
    </h3>
    <div class="search">
     <input class="box" placeholder="Email address" type="text"/>
     <span class="try">
      Try 30 days free
      <i class="fas fa-chevron-right">
      </i>
     </span>
    </div>
    <h4>
    </h4>
   </div>
  </div>
  <div class="container1">
   <div class="text">
    <h1>
    </h1>
    <p>
    </p>
   </div>
   <div class="image">
    <img src="https://thumbs.dreamstime.com/b/netflix-app-tv-screen-playing-series-logo-behind-netflix-app-tv-screen-playing-series-logo-behind-160045137.jpg"/>
   </div>
  </div>
  <div class="container1">
   <div class="image">
    <img src="https://assets.nflxext.com/ffe/siteui/acquisition/ourStory/fuji/desktop/mobile.png"/>
   </div>
   <div class="text">
    <h1>
    </h1>
    <p>
    </p>
   </div>
  </div>
  <div class="container1">
   <div class="text">
    <h1>
    </h1>
    <p>
    </p>
   </div>
   <div class="image">
    <img src="https://assets.nflxext.com/ffe/siteui/acquisition/ourStory/fuji/desktop/device-pile.png"/>
   </div>
  </div>
  <div class="question">
   <h1>
   </h1>
   <div class="quest">
    <div class="textbox">
     What is Netflix?
    </div>
    <i class="las la-plus">
    </i>
   </div>
   <div class="quest">
    <div class="textbox">
     How much does Netflix cost?
    </div>
    <i class="las la-plus">
    </i>
   </div>
   <div class="quest">
    <div class="textbox">
     Where can I watch?
    </div>
    <i class="las la-plus">
    </i>
   </div>
   <div class="quest">
    <div class="textbox">
     How do I cancel?
    </div>
    <i class="las la-plus">
    </i>
   </div>
   <div class="quest">
    <div class="textbox">
     What can I watch on Netflix??
    </div>
    <i class="las la-plus">
    </i>
   </div>
   <div class="quest">
    <div class="textbox">
     What is Netflix?
    </div>
    <i class="las la-plus">
    </i>
   </div>
   <div class="search1">
    <input class="box1" placeholder="Email" type="text"/>
    <span class="try1">
     Try 30 days free
     <i class="fas fa-chevron-right">
     </i>
    </span>
   </div>
   <h4>
   </h4>
  </div>
  <div class="footer">
   <div class="footercon">
    <div class="flex1">
     <h5>
     </h5>
     <h5>
     </h5>
    </div>
    <ul class="list1">
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
    </ul>
    <ul class="list1">
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
    </ul>
    <ul class="list1">
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
    </ul>
    <ul class="list1">
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
     <li>
      <a href="">
      </a>
     </li>
    </ul>
   </div>
  </div>
  <div class="end">
   <h2>
   </h2>
   <h2>
   </h2>
  </div>
 </body>
</html>

"@

$ws.Range("A2").Value = $html
